# Repull data, push all data, mean calculation
# Update column F (dSF) values for rows 2-10 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = -5
